$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 42 (shifts existing rows 42:49 down to 43:50),
# inheriting formatting from row 41 above it.
$ws.Rows.Item(42).Insert()

# Populate the new row 42 with the 2023 TM160 IPA 30 run ("Higher tolls, WFH remains at ~25%")
$ws.Range("A42").Value = 2023
$ws.Range("B42").Value = "2023_TM160_IPA_30"
$ws.Range("C42").Value = "RTP2025_IP"
$ws.Range("D42").Value = "Past year"
$ws.Range("E42").Value = "Higher tolls, WFH remains at ~25%"
$ws.Range("F42").Value = "petrale"
$ws.Range("G42").Value = "n/a"
$ws.Range("H42").Value = "current"
$ws.Range("I42").Value = "BlueprintNetworks_v11\net_2023_Blueprint"
$ws.Range("J42").Value = "model2-b"
$ws.Range("K42").Value = "https://app.asana.com/0/1204085012544660/1205866185692454/f"
$ws.Range("L42").Value = 17.77
$ws.Range("M42").Value = "na"
$ws.Range("N42").Value = "na"
$ws.Range("O42").Value = 0.94
$ws.Range("P42").Value = 0.855
$ws.Range("Q42").Value = 120
$ws.Range("R42").Value = 0
$ws.Range("S42").Value = 45

# Keep the active selection on the newly added row, matching the saved workbook view.
$ws.Range("A42").Select()
